# BB-22361: Impossible to add product to Quick Order Form when multiple
# sellers have the same product SKU.
#
# The quick-order.xlsx functional-test fixture used ambiguous SKU-looking
# values ("1ABSC"/"2ABSC") for its two sample rows; update them to the
# unique "product-N" identifiers used by the updated test data, and leave
# the selection on the last data row as the fixture now does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("1ABSC" -> "product-1") and row 3 ("2ABSC" -> "product-2") hold the
# sample item numbers in column A.
$ws.Range("A2").Value = "product-1"
$ws.Range("A3").Value = "product-2"

# Move the active selection to the last data cell (A3), matching the
# fixture's saved cursor position.
$ws.Range("A3").Select()
